# Amend corrected label annotations
# Update the "labels" column (F) values to match corrected (lower-cased /
# reordered) annotation text per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = "135_product_information"
    9   = "135_product_information"
    12  = "use restrictions || application instructions"
    19  = "93_referral_statement"
    26  = "18_hazards_to_humans_and_domestic_animals"
    27  = "ppe"
    28  = "ppe"
    33  = "off target movement || application instructions || env warning - species || env warning - water"
    40  = "application instructions"
    41  = "135_product_information"
    42  = "application instructions"
    45  = "off target movement"
    46  = "application instructions"
    47  = "application instructions"
    48  = "application instructions"
    49  = "mixing || application instructions"
    50  = "mixing || application instructions"
    51  = "application instructions"
    52  = "irrigation || application instructions || chemigation"
    54  = "irrigation"
    55  = "irrigation"
    56  = "irrigation || chemigation"
    57  = "use restrictions"
    58  = "use restrictions"
    59  = "use restrictions"
    257 = "application instructions"
    258 = "mixing"
    259 = "mixing"
    260 = "mixing"
    261 = "mixing"
    271 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
